# Add {ownerEmail} and {ownerPhoneNumber} placeholder paragraphs (plus a
# trailing blank paragraph) right after the existing {ownerAddress}
# paragraph, matching the author's template update.

$d = $word.ActiveDocument

# Locate the paragraph that contains "{ownerAddress}" - the new content is
# inserted directly after it.
$anchor = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*{ownerAddress}*") {
        $anchor = $para
        break
    }
}

# Create a fresh empty paragraph right after the anchor; its formatting
# (spacing after = 0) is inherited from the anchor paragraph, same as Word
# does when you press Enter at the end of that paragraph.
$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()
$target = $d.Range($newPara.Range.Start, $newPara.Range.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$fragOwnerEmail = "<w:p $ns>" +
    "<w:pPr><w:spacing w:after=`"0`"/></w:pPr>" +
    "<w:r><w:t>{</w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:proofErr w:type=`"gramStart`"/>" +
    "<w:r><w:t>owner</w:t></w:r>" +
    "<w:r><w:t>Email</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r><w:t>}</w:t></w:r>" +
    "</w:p>"

$fragOwnerPhoneNumber = "<w:p $ns>" +
    "<w:pPr><w:spacing w:after=`"0`"/></w:pPr>" +
    "<w:r><w:t>{</w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>owner</w:t></w:r>" +
    "<w:r><w:t>PhoneNumber</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t>}</w:t></w:r>" +
    "</w:p>"

$fragBlank = "<w:p $ns><w:pPr><w:spacing w:after=`"0`"/></w:pPr></w:p>"

$target.InsertXML($fragOwnerEmail + $fragOwnerPhoneNumber + $fragBlank)
